# Remove the "Shopper" role test case row from Sheet1's test data table.
# This was row 14 (test_case_id "singleuseradd_shopper_success", canvas_role
# "Shopper", role_id "aa000g0l"). Deleting the whole row shifts the rows
# below it (Observer, Guest failure) up by one, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(14).Delete()
